$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "MES 4"
$ws.Range("C3").Value = 483148143.2
$ws.Range("D3").Value = 481469708.82
$ws.Range("E3").Value = 9909810.039999999
$ws.Range("F3").Value = 10574124.34
$ws.Range("G3").Value = 493057953.24
$ws.Range("H3").Value = 492043833.16
$ws.Range("I3").Value = 1014120.08
$ws.Range("C4").Value = 190634942.43
$ws.Range("D4").Value = 187382795.66
$ws.Range("E4").Value = 3501504.82
$ws.Range("F4").Value = 3267426.13
$ws.Range("G4").Value = 194136447.25
$ws.Range("H4").Value = 190650221.79
$ws.Range("I4").Value = 3486225.46
$ws.Range("C5").Value = 3522959.97
$ws.Range("D5").Value = 3233940.31
$ws.Range("E5").Value = 68597.13
$ws.Range("F5").Value = 61590.78
$ws.Range("G5").Value = 3591557.1
$ws.Range("H5").Value = 3295531.09
$ws.Range("I5").Value = 296026.01
$ws.Range("D6").Value = 3401113.67
$ws.Range("F6").Value = 692.4
$ws.Range("H6").Value = 3401806.07
$ws.Range("I6").Value = 1805555.45
$ws.Range("C7").Value = 9534611.4
$ws.Range("D7").Value = 10463223.64
$ws.Range("E7").Value = 259972.74
$ws.Range("F7").Value = 243397.5
$ws.Range("G7").Value = 9794584.140000001
$ws.Range("H7").Value = 10706621.14
$ws.Range("J7").Value = 912037
$ws.Range("C8").Value = 39106870.99
$ws.Range("E8").Value = 48349.68
$ws.Range("G8").Value = 39155220.67
$ws.Range("I8").Value = 28027485.26
$ws.Range("C9").Value = 45811120.84
$ws.Range("D9").Value = 41331171.15
$ws.Range("E9").Value = 527917.91
$ws.Range("F9").Value = 601338.03
$ws.Range("G9").Value = 46339038.75
$ws.Range("H9").Value = 41932509.18
$ws.Range("I9").Value = 4406529.57
$ws.Range("D10").Value = 35185144.84
$ws.Range("F10").Value = 634648.34
$ws.Range("H10").Value = 35819793.18
$ws.Range("J10").Value = 25520970.54
$ws.Range("C11").Value = 29489652.82
$ws.Range("D11").Value = 30779593.07
$ws.Range("E11").Value = 903836.05
$ws.Range("F11").Value = 716342.3199999999
$ws.Range("G11").Value = 30393488.87
$ws.Range("H11").Value = 31495935.39
$ws.Range("J11").Value = 1102446.52
$ws.Range("C12").Value = 22416323.67
$ws.Range("D12").Value = 21604592.34
$ws.Range("E12").Value = 309765.46
$ws.Range("F12").Value = 364061.42
$ws.Range("G12").Value = 22726089.13
$ws.Range("H12").Value = 21968653.76
$ws.Range("I12").Value = 757435.37
$ws.Range("C13").Value = 109071433.52
$ws.Range("D13").Value = 110129871.63
$ws.Range("E13").Value = 2393938.11
$ws.Range("F13").Value = 2240790.8
$ws.Range("G13").Value = 111465371.63
$ws.Range("H13").Value = 112370662.43
$ws.Range("J13").Value = 905290.8
$ws.Range("C14").Value = 44805507.88
$ws.Range("D14").Value = 58226524.29
$ws.Range("E14").Value = 700840.13
$ws.Range("F14").Value = 60598.38
$ws.Range("G14").Value = 45506348.01
$ws.Range("H14").Value = 58287122.67
$ws.Range("J14").Value = 12780774.66
$ws.Range("C15").Value = 1925148.37
$ws.Range("D15").Value = 1901416.88
$ws.Range("E15").Value = 31543.56
$ws.Range("F15").Value = 32025.3
$ws.Range("G15").Value = 1956691.93
$ws.Range("H15").Value = 1933442.18
$ws.Range("I15").Value = 23249.75
$ws.Range("C16").Value = 3032007.04
$ws.Range("D16").Value = 77017.57000000001
$ws.Range("E16").Value = 1084112.59
$ws.Range("F16").Value = 32047.73
$ws.Range("G16").Value = 4116119.63
$ws.Range("H16").Value = 109065.3
$ws.Range("I16").Value = 4007054.33
$ws.Range("C17").Value = 654362.75
$ws.Range("D17").Value = 668770.55
$ws.Range("E17").Value = 243967.47
$ws.Range("F17").Value = 259972.74
$ws.Range("G17").Value = 898330.22
$ws.Range("H17").Value = 928743.29
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 30413.07
$ws.Range("C18").Value = 1661343.48
$ws.Range("E18").Value = 446162.19
$ws.Range("G18").Value = 2107505.67
$ws.Range("I18").Value = 2107505.67
$ws.Range("C19").Value = 2027135.58
$ws.Range("E19").Value = 604296.21
$ws.Range("G19").Value = 2631431.79
$ws.Range("I19").Value = 2498267.96
$ws.Range("C20").Value = 79157.74000000001
$ws.Range("E20").Value = 58903.74
$ws.Range("G20").Value = 138061.48
$ws.Range("I20").Value = 138061.48
$ws.Range("C21").Value = 337641.88
$ws.Range("E21").Value = 119142.66
$ws.Range("G21").Value = 456784.54
$ws.Range("I21").Value = 456784.54
$ws.Range("C22").Value = 628545.14
$ws.Range("E22").Value = 140520.69
$ws.Range("G22").Value = 769065.83
$ws.Range("I22").Value = 758504.83
$ws.Range("C23").Value = 2286579.96
$ws.Range("E23").Value = 731192.46
$ws.Range("G23").Value = 3017772.42
$ws.Range("I23").Value = 3017772.42
$ws.Range("D24").Value = 8155756.43
$ws.Range("F24").Value = 2803664.88
$ws.Range("H24").Value = 10959421.31
$ws.Range("J24").Value = 10959421.31
$ws.Range("D25").Value = 52353.38
$ws.Range("F25").Value = 21746.37
$ws.Range("H25").Value = 74099.75
$ws.Range("J25").Value = 74099.75
$ws.Range("D26").Value = 1650795.62
$ws.Range("F26").Value = 121590
$ws.Range("H26").Value = 1772385.62
$ws.Range("J26").Value = 1772385.62
$ws.Range("D27").Value = 205218.93
$ws.Range("F27").Value = 48316.18
$ws.Range("H27").Value = 253535.11
$ws.Range("J27").Value = 253535.11
$ws.Range("D28").Value = 9960985.439999999
$ws.Range("F28").Value = 3136277.53
$ws.Range("H28").Value = 13097262.97
$ws.Range("J28").Value = 12953538.14
$ws.Range("C29").Value = 8382729.65
$ws.Range("E29").Value = 2661954.34
$ws.Range("G29").Value = 11044683.99
$ws.Range("I29").Value = 10911520.16
$ws.Range("C30").Value = 482619.93
$ws.Range("E30").Value = 187585.93
$ws.Range("G30").Value = 670205.86
$ws.Range("I30").Value = 670205.86
$ws.Range("C31").Value = 467090.73
$ws.Range("E31").Value = 146216.58
$ws.Range("G31").Value = 613307.3100000001
$ws.Range("I31").Value = 613307.3100000001
$ws.Range("C32").Value = 628545.13
$ws.Range("E32").Value = 140520.68
$ws.Range("G32").Value = 769065.8100000001
$ws.Range("I32").Value = 758504.8100000001
$ws.Range("G33").Value = 1041005034.26
$ws.Range("H33").Value = 1042515830.46
$ws.Range("I33").Value = 65754116.32
$ws.Range("J33").Value = 67264912.52
